$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1586067701024982
$ws.Range("C2").Value = 1.518426349632036
$ws.Range("D2").Value = 7.036628225930749
$ws.Range("E2").Value = 2.652664363603272
$ws.Range("F2").Value = 2.713308999050573
